$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = -17.20414195516174
$ws.Range("J25").Value = -16.98434150865444
$ws.Range("K25").Value = -3.386369312788844
$ws.Range("I26").Value = -17.22295426235204
$ws.Range("J26").Value = -3.624982066486445
$ws.Range("K26").Value = -0.9294641821589158
$ws.Range("H27").Value = -17.25689497300509
$ws.Range("I27").Value = -3.658922777139496
$ws.Range("J27").Value = -0.9634048928119666
$ws.Range("K27").Value = 0.3458441452769406
$ws.Range("G28").Value = -17.13529955079305
$ws.Range("H28").Value = -3.537327354927458
$ws.Range("I28").Value = -0.8418094705999279
$ws.Range("J28").Value = 0.4674395674889793
$ws.Range("K28").Value = 1.988799384936469
$ws.Range("F29").Value = -17.39412186950338
$ws.Range("G29").Value = -3.79614967363779
$ws.Range("H29").Value = -1.10063178931026
$ws.Range("I29").Value = 0.2086172487786471
$ws.Range("J29").Value = 1.729977066226137
$ws.Range("K29").Value = -4.035764575904229
$ws.Range("E30").Value = -17.46219349820464
$ws.Range("F30").Value = -3.864221302339042
$ws.Range("G30").Value = -1.168703418011513
$ws.Range("H30").Value = 0.140545620077394
$ws.Range("I30").Value = 1.661905437524884
$ws.Range("J30").Value = -4.103836204605481
$ws.Range("K30").Value = 2.803854532616474
$ws.Range("D31").Value = -17.61617421205837
$ws.Range("E31").Value = -4.018202016192776
$ws.Range("F31").Value = -1.322684131865246
$ws.Range("G31").Value = -0.01343509377633867
$ws.Range("H31").Value = 1.507924723671151
$ws.Range("I31").Value = -4.257816918459215
$ws.Range("J31").Value = 2.649873818762741
$ws.Range("K31").Value = 0.2395698154090965
$ws.Range("C32").Value = -18.09929231679948
$ws.Range("D32").Value = -4.501320120933885
$ws.Range("E32").Value = -1.805802236606356
$ws.Range("F32").Value = -0.4965531985174484
$ws.Range("G32").Value = 1.024806618930041
$ws.Range("H32").Value = -4.740935023200324
$ws.Range("I32").Value = 2.166755714021631
$ws.Range("J32").Value = -0.2435482893320133
$ws.Range("K32").Value = -1.938804334703323
$ws.Range("B33").Value = -20.4108400473813
$ws.Range("C33").Value = -6.812867851515707
$ws.Range("D33").Value = -4.117349967188177
$ws.Range("E33").Value = -2.80810092909927
$ws.Range("F33").Value = -1.286741111651781
$ws.Range("G33").Value = -7.052482753782146
$ws.Range("H33").Value = -0.1447920165601908
$ws.Range("I33").Value = -2.555096019913835
$ws.Range("J33").Value = -4.250352065285145
$ws.Range("K33").Value = -1.311034623099504
$ws.Range("B34").Value = -4.980277842704087
$ws.Range("C34").Value = -2.284759958376557
$ws.Range("D34").Value = -0.9755109202876501
$ws.Range("E34").Value = 0.5458488971598395
$ws.Range("F34").Value = -5.219892744970526
$ws.Range("G34").Value = 1.68779799225143
$ws.Range("H34").Value = -0.7225060111022149
$ws.Range("I34").Value = -2.417762056473524
$ws.Range("J34").Value = 0.5215553857121161
$ws.Range("K34").Value = -0.8573456354163971
$ws.Range("B35").Value = -1.909452872482039
$ws.Range("C35").Value = -0.6002038343931317
$ws.Range("D35").Value = 0.9211559830543579
$ws.Range("E35").Value = -4.844585659076007
$ws.Range("F35").Value = 2.063105078145948
$ws.Range("G35").Value = -0.3471989252076966
$ws.Range("H35").Value = -2.042454970579006
$ws.Range("I35").Value = 0.8968624716066345
$ws.Range("J35").Value = -0.4820385495218787
$ws.Range("K35").Value = 1.022680634228276
$ws.Range("B36").Value = -0.3749684946957029
$ws.Range("C36").Value = 1.146391322751787
$ws.Range("D36").Value = -4.619350319378579
$ws.Range("E36").Value = 2.288340417843377
$ws.Range("F36").Value = -0.1219635855102677
$ws.Range("G36").Value = -1.817219630881577
$ws.Range("H36").Value = 1.122097811304063
$ws.Range("I36").Value = -0.2568032098244498
$ws.Range("J36").Value = 1.247915973925705
$ws.Range("K36").Value = 0.9211944755864938
$ws.Range("B37").Value = 1.14167028642729
$ws.Range("C37").Value = -4.624071355703076
$ws.Range("D37").Value = 2.283619381518879
$ws.Range("E37").Value = -0.1266846218347649
$ws.Range("F37").Value = -1.821940667206074
$ws.Range("G37").Value = 1.117376774979566
$ws.Range("H37").Value = -0.261524246148947
$ws.Range("I37").Value = 1.243194937601208
$ws.Range("J37").Value = 0.9164734392619965
$ws.Range("K37").Value = 1.754572496573351
$ws.Range("B38").Value = -4.896022371537698
$ws.Range("C38").Value = 2.011668365684257
$ws.Range("D38").Value = -0.3986356376693871
$ws.Range("E38").Value = -2.093891683040697
$ws.Range("F38").Value = 0.8454257591449439
$ws.Range("G38").Value = -0.5334752619835692
$ws.Range("H38").Value = 0.9712439217665854
$ws.Range("I38").Value = 0.6445224234273743
$ws.Range("J38").Value = 1.482621480738728
$ws.Range("K38").Value = -1.197078432822523
$ws.Range("B39").Value = 2.832994207660627
$ws.Range("C39").Value = 0.4226902043069828
$ws.Range("D39").Value = -1.272565841064327
$ws.Range("E39").Value = 1.666751601121314
$ws.Range("F39").Value = 0.2878505799928007
$ws.Range("G39").Value = 1.792569763742955
$ws.Range("H39").Value = 1.465848265403744
$ws.Range("I39").Value = 2.303947322715098
$ws.Range("J39").Value = -0.3757525908461526
$ws.Range("K39").Value = 0.8596701032167943
$ws.Range("B40").Value = -0.03790361708925488
$ws.Range("C40").Value = -1.733159662460564
$ws.Range("D40").Value = 1.206157779725076
$ws.Range("E40").Value = -0.172743241403437
$ws.Range("F40").Value = 1.331975942346718
$ws.Range("G40").Value = 1.005254444007507
$ws.Range("H40").Value = 1.843353501318861
$ws.Range("I40").Value = -0.8363464122423903
$ws.Range("J40").Value = 0.3990762818205566
$ws.Range("K40").Value = 2.123380570563001
$ws.Range("B41").Value = -1.502432366452369
$ws.Range("C41").Value = 1.436885075733271
$ws.Range("D41").Value = 0.05798405460475808
$ws.Range("E41").Value = 1.562703238354913
$ws.Range("F41").Value = 1.235981740015702
$ws.Range("G41").Value = 2.074080797327056
$ws.Range("H41").Value = -0.6056191162341953
$ws.Range("I41").Value = 0.6298035778287516
$ws.Range("J41").Value = 2.354107866571197
$ws.Range("K41").Value = 2.720553059184225
$ws.Range("B42").Value = 2.021185630531559
$ws.Range("C42").Value = 0.6422846094030465
$ws.Range("D42").Value = 2.147003793153201
$ws.Range("E42").Value = 1.82028229481399
$ws.Range("F42").Value = 2.658381352125344
$ws.Range("G42").Value = -0.02131856143590682
$ws.Range("H42").Value = 1.21410413262704
$ws.Range("I42").Value = 2.938408421369485
$ws.Range("J42").Value = 3.304853613982513
$ws.Range("K42").Value = -2.170593826049543
$ws.Range("B43").Value = 2.59639870328499
$ws.Range("C43").Value = 4.101117887035145
$ws.Range("D43").Value = 3.774396388695934
$ws.Range("E43").Value = 4.612495446007288
$ws.Range("F43").Value = 1.932795532446037
$ws.Range("G43").Value = 3.168218226508984
$ws.Range("H43").Value = 4.892522515251429
$ws.Range("I43").Value = 5.258967707864457
$ws.Range("J43").Value = -0.2164797321675991
$ws.Range("K43").Value = 3.144233349489796
$ws.Range("B44").Value = 2.341506873006513
$ws.Range("C44").Value = 2.014785374667302
$ws.Range("D44").Value = 2.852884431978656
$ws.Range("E44").Value = 0.173184518417405
$ws.Range("F44").Value = 1.408607212480352
$ws.Range("G44").Value = 3.132911501222797
$ws.Range("H44").Value = 3.499356693835825
$ws.Range("I44").Value = -1.976090746196231
$ws.Range("J44").Value = 1.384622335461164
$ws.Range("B45").Value = 1.466561183616531
$ws.Range("C45").Value = 2.304660240927885
$ws.Range("D45").Value = -0.3750396726333658
$ws.Range("E45").Value = 0.8603830214295811
$ws.Range("F45").Value = 2.584687310172026
$ws.Range("G45").Value = 2.951132502785054
$ws.Range("H45").Value = -2.524314937247002
$ws.Range("I45").Value = 0.8363981444103927
$ws.Range("B46").Value = 2.012044252255831
$ws.Range("C46").Value = -0.6676556613054199
$ws.Range("D46").Value = 0.567767032757527
$ws.Range("E46").Value = 2.292071321499972
$ws.Range("F46").Value = 2.658516514113
$ws.Range("G46").Value = -2.816930925919056
$ws.Range("H46").Value = 0.5437821557383387
$ws.Range("B47").Value = -0.9880269291625079
$ws.Range("C47").Value = 0.247395764900439
$ws.Range("D47").Value = 1.971700053642884
$ws.Range("E47").Value = 2.338145246255912
$ws.Range("F47").Value = -3.137302193776144
$ws.Range("G47").Value = 0.2234108878812506
$ws.Range("B48").Value = 0.5533907096283328
$ws.Range("C48").Value = 2.277694998370778
$ws.Range("D48").Value = 2.644140190983806
$ws.Range("E48").Value = -2.83130724904825
$ws.Range("F48").Value = 0.5294058326091444
$ws.Range("B49").Value = 1.9047312492914
$ws.Range("C49").Value = 2.271176441904428
$ws.Range("D49").Value = -3.204270998127628
$ws.Range("E49").Value = 0.1564420835297668
$ws.Range("B50").Value = 2.102905296315023
$ws.Range("C50").Value = -3.372542143717033
$ws.Range("D50").Value = -0.01182906205963841
$ws.Range("B51").Value = -3.650852129092033
$ws.Range("C51").Value = -0.2901390474346385
$ws.Range("B52").Value = -0.2631055417942008
